$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the magnet SKU from D73-N52 to D83-N52
$ws.Range("B4").Value = "D83-N52"
$ws.Range("C4").Value = "D83-N52"

# Update the price formula for row 4 to a fixed hardcoded value
$ws.Range("E4").Formula = "=1.71"

# Update row 12's price to reference E4 instead of a hardcoded value
$ws.Range("E12").Formula = "=E4"

# Update selection to E13 (matches final cursor position in diff)
$ws.Range("E13").Select()
